$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 218.721583
$ws.Range("N2").Value = 656.164749
$ws.Range("O2").Value = 0.7793342808141792
$ws.Range("P2").Value = 0.7793342808141792
$ws.Range("Q2").Value = 73.50867868658332
$ws.Range("R2").Value = 661.57810817925
$ws.Range("S2").Value = 0.7793342808141792
$ws.Range("T2").Value = 0.7793342808141792

$ws.Range("O3").Value = 0.164954193449581
$ws.Range("P3").Value = 0.164954193449581
$ws.Range("S3").Value = 0.164954193449581
$ws.Range("T3").Value = 0.164954193449581

$ws.Range("M4").Value = 8.515309999999999
$ws.Range("N4").Value = 25.54593
$ws.Range("O4").Value = 0.03034118948727519
$ws.Range("P4").Value = 0.03034118948727519
$ws.Range("Q4").Value = 2.861853769166666
$ws.Range("R4").Value = 25.7566839225
$ws.Range("S4").Value = 0.03034118948727519
$ws.Range("T4").Value = 0.03034118948727519

$ws.Range("M5").Value = 7.120231
$ws.Range("N5").Value = 21.360693
$ws.Range("O5").Value = 0.02537033624896462
$ws.Range("P5").Value = 0.02537033624896462
$ws.Range("Q5").Value = 2.392990968583333
$ws.Range("R5").Value = 21.53691871725
$ws.Range("S5").Value = 0.02537033624896462
$ws.Range("T5").Value = 0.02537033624896462
